# TaiGer Transcript-Program_Comparer / CS_Programs.xlsx
# Populate the "Program_choosing" sheet with the full list of CS programs
# and extend the Yes/No data validation + selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the order in which new strings are first written controls the order
# they land in in the shared-strings table, so these assignments are
# intentionally not in simple top-to-bottom row order.
$ws.Range("A3").Value  = "RWTH_Aachen_Data Science"
$ws.Range("B3").Value  = "Yes"

$ws.Range("A4").Value  = "Freie Uni Berlin - Data Science"
$ws.Range("B4").Value  = "Yes"

$ws.Range("A5").Value  = "TU Berlin Computer Science"
$ws.Range("B5").Value  = "Yes"

# Row 2 used to hold "TUM_Info" / "Yes" - rename the program to its full name.
$ws.Range("A2").Value  = "TUM_Informatics"
$ws.Range("B2").Value  = "Yes"

$ws.Range("A6").Value  = "TUM Data Engineering and Analytics"
$ws.Range("B6").Value  = "Yes"

$ws.Range("A7").Value  = "TU Delft Computer Science"
$ws.Range("B7").Value  = "Yes"

$ws.Range("A8").Value  = "RWTH_Aachen_DDS"
$ws.Range("B8").Value  = "Yes"

$ws.Range("A9").Value  = "RWTH_Aachen_TIME"
$ws.Range("B9").Value  = "Yes"

$ws.Range("A10").Value = "Uni_Goettingen_Applied_CS"
$ws.Range("B10").Value = "Yes"

# Extend the Yes/No list validation from B1:B5 to cover the new rows, B1:B10.
$ws.Range("B1:B5").Validation.Delete()
$ws.Range("B1:B10").Validation.Add(3, 1, 1, '"Yes,No"')

# Leave the selection where the user would land after typing the last row.
$ws.Range("A11").Select() | Out-Null
